$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 641.3333
$ws.Range("I38").Value = 169.6
$ws.Range("K38").Value = 508.8
$ws.Range("M38").Value = -136.8

# Row 39
$ws.Range("H39").Value = 148.6
$ws.Range("I39").Value = 94.083336
$ws.Range("J39").Value = 366.66666
$ws.Range("K39").Value = 282.250008
$ws.Range("L39").Value = 1099.99998
$ws.Range("M39").Value = 13.74999200000002
$ws.Range("N39").Value = -1691.99998

# Row 134
$ws.Range("H134").Value = 45000
$ws.Range("J134").Value = 45000
$ws.Range("L134").Value = 45000
$ws.Range("N134").Value = -55140

# Row 135
$ws.Range("H135").Value = 34830.965
$ws.Range("I135").Value = 34830.965
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 313478.6849999999
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -310943.6849999999
$ws.Range("N135").ClearContents()

# Row 137
$ws.Range("H137").Value = 3128537
$ws.Range("I137").Value = 4169126.2
$ws.Range("J137").Value = 6769.125
$ws.Range("K137").Value = 12507378.6
$ws.Range("L137").Value = 20307.375
$ws.Range("M137").Value = -12504828.6
$ws.Range("N137").Value = -25407.375

# Row 138
$ws.Range("H138").Value = 2647.66
$ws.Range("I138").Value = 1901.4
$ws.Range("J138").Value = 2967.4856
$ws.Range("K138").Value = 5704.200000000001
$ws.Range("L138").Value = 8902.4568
$ws.Range("M138").Value = -564.2000000000007
$ws.Range("N138").Value = -19182.4568

# Row 139
$ws.Range("H139").Value = 40909.09
$ws.Range("J139").Value = 40909.09
$ws.Range("L139").Value = 40909.09
$ws.Range("N139").Value = -51189.09

# Row 141
$ws.Range("H141").Value = 1330.9783
$ws.Range("I141").Value = 1100.4762
$ws.Range("J141").Value = 3751.25
$ws.Range("K141").Value = 3301.4286
$ws.Range("L141").Value = 11253.75
$ws.Range("M141").Value = 1878.5714
$ws.Range("N141").Value = -21613.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 50101476
$ws.Range("I61").Value = 66734900
$ws.Range("J61").Value = 201203.8
$ws.Range("K61").Value = 66734900
$ws.Range("L61").Value = 201203.8
$ws.Range("M61").Value = -66734688
$ws.Range("N61").Value = -201627.8

# Row 101
$ws.Range("H101").Value = 40000
$ws.Range("J101").Value = 40000
$ws.Range("L101").Value = 40000
$ws.Range("N101").Value = -46490

# Row 102
$ws.Range("H102").Value = 10991707
$ws.Range("I102").Value = 14288410
$ws.Range("J102").Value = 2700
$ws.Range("K102").Value = 14288410
$ws.Range("L102").Value = 2700
$ws.Range("M102").Value = -14286788
$ws.Range("N102").Value = -5944

# Row 136
$ws.Range("H136").Value = 50101476
$ws.Range("I136").Value = 66734900
$ws.Range("J136").Value = 201203.8
$ws.Range("K136").Value = 200204700
$ws.Range("L136").Value = 603611.3999999999
$ws.Range("M136").Value = -200202150
$ws.Range("N136").Value = -608711.3999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 60
$ws.Range("H60").Value = 49800
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 49800
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 49800
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -50998

# Row 135
$ws.Range("H135").Value = 57911.8
$ws.Range("J135").Value = 57911.8
$ws.Range("L135").Value = 57911.8
$ws.Range("N135").Value = -68051.8

# Row 138
$ws.Range("H138").Value = 32842.223
$ws.Range("J138").Value = 32842.223
$ws.Range("L138").Value = 32842.223
$ws.Range("N138").Value = -43122.223

# Row 140
$ws.Range("H140").Value = 71347.5
$ws.Range("J140").Value = 71347.5
$ws.Range("L140").Value = 71347.5
$ws.Range("N140").Value = -81707.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 1236
$ws.Range("I99").Value = 942.2
$ws.Range("J99").Value = 1603.25
$ws.Range("K99").Value = 942.2
$ws.Range("L99").Value = 1603.25
$ws.Range("M99").Value = 555.8
$ws.Range("N99").Value = -4599.25

# Row 122
$ws.Range("H122").Value = 2219.56
$ws.Range("I122").Value = 1703.6923
$ws.Range("J122").Value = 2778.4167
$ws.Range("K122").Value = 5111.0769
$ws.Range("L122").Value = 8335.250100000001
$ws.Range("M122").Value = -2661.0769
$ws.Range("N122").Value = -13235.2501

# Row 126
$ws.Range("H126").Value = 1236
$ws.Range("I126").Value = 942.2
$ws.Range("J126").Value = 1603.25
$ws.Range("K126").Value = 2826.6
$ws.Range("L126").Value = 4809.75
$ws.Range("M126").Value = -356.6000000000004
$ws.Range("N126").Value = -9749.75

# Row 132
$ws.Range("H132").Value = 21240.686
$ws.Range("I132").Value = 1327.6154
$ws.Range("J132").Value = 85958.164
$ws.Range("K132").Value = 3982.8462
$ws.Range("L132").Value = 257874.492
$ws.Range("M132").Value = -1452.8462
$ws.Range("N132").Value = -262934.492

# Row 135
$ws.Range("H135").Value = 49290
$ws.Range("J135").Value = 49290
$ws.Range("L135").Value = 49290
$ws.Range("N135").Value = -59430

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 863.84375
$ws.Range("J107").Value = 1032.9
$ws.Range("L107").Value = 3098.7
$ws.Range("N107").Value = -6938.700000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 75
$ws.Range("H75").Value = 46960
$ws.Range("J75").Value = 46960
$ws.Range("L75").Value = 46960
$ws.Range("N75").Value = -48708

# Row 78
$ws.Range("H78").Value = 46960
$ws.Range("J78").Value = 46960
$ws.Range("L78").Value = 140880
$ws.Range("N78").Value = -149616

# Row 97
$ws.Range("H97").Value = 3108.75
$ws.Range("I97").Value = 3081.4285
$ws.Range("J97").Value = 3300
$ws.Range("K97").Value = 3081.4285
$ws.Range("L97").Value = 3300
$ws.Range("M97").Value = -2585.4285
$ws.Range("N97").Value = -4292

# Row 133
$ws.Range("H133").Value = 61960
$ws.Range("J133").Value = 61960
$ws.Range("L133").Value = 61960
$ws.Range("N133").Value = -72080

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1594.3
$ws.Range("I40").Value = 1498
$ws.Range("J40").Value = 1738.75
$ws.Range("K40").Value = 1498
$ws.Range("L40").Value = 1738.75
$ws.Range("M40").Value = -1362
$ws.Range("N40").Value = -2010.75

# Row 100
$ws.Range("H100").Value = 1684.3684
$ws.Range("I100").Value = 1250.375
$ws.Range("K100").Value = 1250.375
$ws.Range("M100").Value = -709.375

# Row 132
$ws.Range("H132").Value = 26636.28
$ws.Range("I132").Value = 2480.4849
$ws.Range("J132").Value = 106350.4
$ws.Range("K132").Value = 7441.4547
$ws.Range("L132").Value = 319051.2
$ws.Range("M132").Value = -4911.4547
$ws.Range("N132").Value = -324111.2

# Row 134
$ws.Range("H134").Value = 43022.062
$ws.Range("J134").Value = 43022.062
$ws.Range("L134").Value = 43022.062
$ws.Range("N134").Value = -53162.062

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 568.0833
$ws.Range("I126").Value = 568.0833
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 1704.2499
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 765.7501
$ws.Range("N126").ClearContents()
